$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New job listing data for rows 2-7 (header row 1 unchanged)
$data = @(
    @("Entry Level Front End Developer", "Revature", "Brockton, MA On-site 2 hours ago", "11 applicants", "https://revature.com/jobs/entry-level-software-engineer-3/?utm_source=linkedin&sourcedBy=BalaLP"),
    @("Front End Developer - Freelance [Remote]", "Braintrust", "Boston, MA Remote 1 day ago", "112 applicants", "https://boards.greenhouse.io/braintrust"),
    @("Front End Engineer", "GavinHeath", "United States Remote 2 hours ago", "70 applicants", "https://www1.jobdiva.com/portal/?a=k1jdnw97y2ba69mixpqijkgh9633e3075ay54ctusnobplysx8gv35gw2tu44062&jobid=24432843#/jobs/24432843?compid=0&SearchString=&StatesString=&id=24432843&source=linkedin.com"),
    @("Front-End Developer", "TekVivid, Inc", "Tewksbury, MA Hybrid 1 hour ago", "6 applicants", "https://tekvivid.oorwin.com/careers/index.html#/job/e26973e6ee8ab9cd8cb3f207d1b90f00d2669eff?source=linkedin_free&postedby=16"),
    @("Software Engineer (Front-End)-Minneapolis, MN: 23-01822", "Akraya, Inc.", "Boston, MA 12 hours ago", "7 applicants", "https://www2.jobdiva.com/portal/?a=vmjdnwvuqk8ziib1dp1dwwtj5lmlhe03a7vhen21vtdgmiqsu14lwy7xz479eb4w&jobid=20120600#/jobs/20120600?compid=0&SearchString=&StatesString=&id=20120600&source=linkedin.com"),
    @("Store Associate, PT - Boston", "Nike", "Boston, MA 1 day ago", "2 applicants", "https://jobs.nike.com/job/R-136?source=BY_LinkedIn_OrganicFeed")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $rowData = $data[$i]
    $ws.Cells.Item($row, 1).Value = $rowData[0]
    $ws.Cells.Item($row, 2).Value = $rowData[1]
    $ws.Cells.Item($row, 3).Value = $rowData[2]
    $ws.Cells.Item($row, 4).Value = $rowData[3]
    $ws.Cells.Item($row, 5).Value = $rowData[4]
}

# Remove now-unused rows 8-12 (sheet shrinks from A1:E12 to A1:E7)
$ws.Range("A8:E12").Delete() | Out-Null
